# I-No picks DP problems: add 3 new Dynamic Programming rows to the
# "Easy" sheet (rows 64-66), matching the formatting already used by the
# other Dynamic Programming rows above (61-63).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{ Row = 64; Num = 464; Question = "Can I Win" },
    @{ Row = 65; Num = 322; Question = "Coin Change" },
    @{ Row = 66; Num = 152; Question = "Maximum Product Subarray" }
)

foreach ($r in $newRows) {
    $row = $r.Row

    # Copy the formatting from the last existing DP row (63) down onto the
    # new row for the columns that carry a non-default style: #, Question,
    # PIC, Difficulty, Language. (Status/Discussed/Comment are left blank
    # for these new entries, same as the source edit.)
    $ws.Range("A63:B63").Copy() | Out-Null
    $ws.Range("A${row}:B${row}").PasteSpecial(-4122) | Out-Null

    $ws.Range("D63").Copy() | Out-Null
    $ws.Range("D${row}").PasteSpecial(-4122) | Out-Null

    $ws.Range("F63:G63").Copy() | Out-Null
    $ws.Range("F${row}:G${row}").PasteSpecial(-4122) | Out-Null

    $ws.Cells.Item($row, 1).Value = $r.Num          # # (LeetCode problem number)
    $ws.Cells.Item($row, 2).Value = $r.Question      # Question
    $ws.Cells.Item($row, 3).Value = "Dynamic Programming"  # Topic
    $ws.Cells.Item($row, 4).Value = "I-No"           # PIC
    $ws.Cells.Item($row, 6).Value = "Medium"         # Difficulty
    $ws.Cells.Item($row, 7).Value = "Python"         # Language
}

$excel.ActiveWindow.ScrollRow = 49
$ws.Range("C71").Select() | Out-Null
